$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$chartObj = $ws.Shapes.AddChart2(227, 5)
$chart = $chartObj.Chart
$chart.SetSourceData($ws.Range("D10:E14"))
$chart.SeriesCollection(1).ApplyDataLabels(1, $false, $false, $false, $false, $false, $true, $true, $false)
Write-Host "done"
